# Fruta / hortaliza, semanal
# Re-shuffles the per-record block (Fecha, Calidad, Volumen, Precio min/max/prom,
# Unidad de comercializacion, Origen, Precio $/Kg, Kg/unidad) across the
# existing data rows (rows 2-24) of the active sheet. Row identity columns
# (Mercado ID, Mercado, Region, Codreg, Tipo, Producto*, Categoria*, Variedad)
# stay put; only the "event" columns D and L..T move between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that make up the record "payload" that gets reshuffled.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot the current (pre-edit) values for every payload column on every
# data row BEFORE any writes happen, since several rows swap with each other.
$snapshot = @{}
foreach ($row in 2..24) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowVals
}

# Target row -> source row (i.e. target row receives the payload that used
# to live on the source row).
$mapping = @{
    2  = 18
    3  = 6
    4  = 7
    5  = 5
    6  = 3
    7  = 22
    8  = 16
    9  = 15
    10 = 2
    11 = 4
    12 = 17
    13 = 24
    14 = 13
    15 = 10
    16 = 19
    17 = 14
    18 = 20
    19 = 11
    20 = 12
    21 = 23
    22 = 21
    23 = 8
    24 = 9
}

foreach ($targetRow in 2..24) {
    $sourceRow = $mapping[$targetRow]
    $src = $snapshot[$sourceRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value = $src[$col]
    }
}
